$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.64129999999998
$ws.Range("A6").Value = -22.92820000000001
$ws.Range("E6").Value = 16.7187
$ws.Range("A7").Value = -21.9333
$ws.Range("B7").Value = 4.9461
$ws.Range("E8").Value = 16.3246
$ws.Range("B12").Value = 4.722699999999999
$ws.Range("D12").Value = -7.158299999999998
$ws.Range("C13").Value = -13.98789999999998
$ws.Range("C14").Value = -14.01459999999999
$ws.Range("B15").Value = 5.141499999999995
$ws.Range("A16").Value = -21.68719999999998
$ws.Range("C16").Value = -11.91249999999999
$ws.Range("E18").Value = 18.48120000000002
$ws.Range("C19").Value = -13.07200000000001
$ws.Range("E19").Value = 16.2797
$ws.Range("A20").Value = -22.72980000000002
$ws.Range("B20").Value = 5.101599999999995
$ws.Range("E20").Value = 16.5513
$ws.Range("B21").Value = 10.66470000000001
$ws.Range("E21").Value = 16.68160000000001
$ws.Range("B22").Value = 9.230100000000002
$ws.Range("C22").Value = -12.76620000000001
$ws.Range("D22").Value = -8.232100000000008
$ws.Range("B23").Value = 9.2224
$ws.Range("E23").Value = 16.0348
$ws.Range("E24").Value = 16.74390000000001
$ws.Range("A28").Value = -21.99729999999999
$ws.Range("A29").Value = -21.844
$ws.Range("B29").Value = 5.084200000000003
$ws.Range("D29").Value = -7.886100000000003
$ws.Range("A32").Value = -21.33580000000001
$ws.Range("B34").Value = 9.009700000000004
$ws.Range("D34").Value = -7.688700000000004
$ws.Range("E35").Value = 15.6084
$ws.Range("C36").Value = -12.61570000000001
$ws.Range("E37").Value = 16.2765
$ws.Range("E39").Value = 15.7346
$ws.Range("A40").Value = -19.25639999999999
$ws.Range("E41").Value = 16.10939999999999
$ws.Range("B42").Value = 10.19059999999999
$ws.Range("B43").Value = 5.6979
$ws.Range("D43").Value = -8.170900000000003
$ws.Range("B44").Value = 4.947200000000002
$ws.Range("B45").Value = 5.369100000000002
$ws.Range("A46").Value = -22.17259999999999
$ws.Range("B46").Value = 5.6549
$ws.Range("C46").Value = -12.30759999999999
$ws.Range("E46").Value = 16.5964
$ws.Range("D48").Value = -7.7659
$ws.Range("B50").Value = 4.973799999999997
$ws.Range("C50").Value = -13.85649999999999
$ws.Range("A51").Value = -22.20810000000001
$ws.Range("B51").Value = 5.662699999999997
$ws.Range("A52").Value = -22.24100000000001
$ws.Range("A57").Value = -22.70080000000001
$ws.Range("E58").Value = 16.45660000000001
$ws.Range("A59").Value = -21.97610000000001
$ws.Range("D60").Value = -8.363199999999999
$ws.Range("E60").Value = 16.3084
$ws.Range("A62").Value = -22.16890000000001
$ws.Range("A66").Value = -21.53129999999997
$ws.Range("B66").Value = 4.873199999999997
$ws.Range("B67").Value = 4.979899999999998
$ws.Range("D68").Value = -7.073199999999995
$ws.Range("D70").Value = -7.7255
$ws.Range("A73").Value = -19.79730000000001
$ws.Range("D73").Value = -7.926100000000003
$ws.Range("E73").Value = 16.3747
$ws.Range("A74").Value = -21.96849999999998
$ws.Range("E76").Value = 16.38450000000001
$ws.Range("B79").Value = 9.83180000000001
$ws.Range("B84").Value = 5.503000000000001
$ws.Range("E85").Value = 16.4991
$ws.Range("D87").Value = -8.741799999999996
$ws.Range("A92").Value = -21.55530000000002
$ws.Range("B92").Value = 4.615099999999996
$ws.Range("D92").Value = -6.189099999999997
$ws.Range("C95").Value = -11.6849
$ws.Range("B97").Value = 5.8777
$ws.Range("C97").Value = -12.0825
$ws.Range("E98").Value = 15.6829
$ws.Range("A100").Value = -22.01309999999999
$ws.Range("D101").Value = -7.871900000000003
